# UsersContext.xlsx — "Add files via upload"
#
# The workbook's Sheet1 holds a flat export of employee records with
# boolean (0/1) flag columns Is_DCEO (Q), Is_Manger (R), Is_TeamLeader (S).
# This re-upload corrects three of those flags:
#   - Row 4 (Marwan, Ahmad...): Is_Manger flips 0 -> 1
#   - Row 9 (the long-titled record): Is_DCEO flips 0 -> 1 and
#     Is_TeamLeader flips 1 -> 0 (the "is a team leader" flag was moved to
#     "is DCEO" for that person).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("R4").Value = 1

$ws.Range("Q9").Value = 1
$ws.Range("S9").Value = 0

# Cosmetic: the re-saved file also shows the author scrolled/selected a
# different cell before uploading.
$ws.Range("M7").Select() | Out-Null
